# Apply the edit described in the commit: update a few header/label cells,
# remove the data that used to live in column G, and move the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update relabeled cells ---
$ws.Range("B1").Value = "Company Service"
$ws.Range("C1").Value = "Use of the Data"
$ws.Range("A2").Value = "Negative User Impact"
$ws.Range("A3").Value = "User Data Confidentiality Violation "

# --- Clear out column G (G1:G6), which is no longer used ---
$ws.Range("G1:G6").ClearContents()

# --- Move the active selection to A17 (matches the saved selection state) ---
$ws.Range("A17").Select()
